{"js": "// Update the MELD / CNN / Wav2Vec result cell: \"0.481/0.171\" -> \"0.481/0.149\"\nconst oldText = \"0.481/0.171\";\nconst newText = \"0.481/0.149\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find text \"${oldText}\" in the document.`);\n}\n\n// Replace every occurrence found (expected to be exactly one) in place,\n// preserving the surrounding run formatting (e.g. lang=\"pl-PL\").\nfor (const range of results.items) {\n  range.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the MELD / CNN / Wav2Vec result cell: \"0.481/0.171\" -> \"0.481/0.149\"\n\n# Word constants (not predefined by the host, so declare them explicitly)\n$wdFindContinue = 1\n$wdReplaceAll    = 2\n\n$oldText = \"0.481/0.171\"\n$newText = \"0.481/0.149\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    $oldText,        # FindText\n    $true,           # MatchCase\n    $false,          # MatchWholeWord\n    $false,          # MatchWildcards\n    $false,          # MatchSoundsLike\n    $false,          # MatchAllWordForms\n    $true,           # Forward\n    $wdFindContinue, # Wrap\n    $false,          # Format\n    $newText,        # ReplaceWith\n    $wdReplaceAll    # Replace\n)\n\nif (-not $found) {\n    throw \"Could not find text '$oldText' in the document.\"\n}\n\nWrite-Output \"Replaced '$oldText' with '$newText': $found\"\n"}
